# Generate Report for Handback
#
# Updates the handback-status workbook so that the two tracked e2e
# markdown files are reported under their new (regenerated) names:
#   62288978-c74f-438f-83ff-b02031d3c663  ->  8cdfdab3-0d92-4118-9672-4add81af50ff
#   ad84fe27-7571-4181-b764-17442849d730  ->  ffffec629dde-426b-428a-a92c-5c9ae518e34f
# together with refreshed xliff hashes and handoff/handback timestamps.
#
# NOTE: Range.Hyperlinks.Delete() clears *every* hyperlink on the sheet
# (not just the target range) in this COM host, and Hyperlinks.Add()
# appends to the collection without de-duplicating an existing entry on
# the same cell. So for each worksheet we delete the whole collection
# exactly once and then re-add every hyperlinked cell (changed or not)
# in a single pass.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2 (file #1)
$ov.Range("A2").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.md"
$ov.Range("G2").Value = "2016-08-13 19:19:38"

# Row 3 (file #2)
$ov.Range("A3").Value = "ffffec629dde-426b-428a-a92c-5c9ae518e34f.md"
$ov.Range("G3").Value = "2016-08-13 19:19:38"

# Rebuild both hyperlinks on this sheet (B2, B3) in one pass.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/62288978-c74f-438f-83ff-b02031d3c663.md", "", "", "e2e\8cdfdab3-0d92-4118-9672-4add81af50ff.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/ad84fe27-7571-4181-b764-17442849d730.md", "", "", "e2e\ffffec629dde-426b-428a-a92c-5c9ae518e34f.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2 (file #1)
$zh.Range("G2").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-13 19:19:30"
$zh.Range("J2").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-13 19:19:59"

# Row 3 (file #2)
$zh.Range("G3").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-13 19:19:30"
$zh.Range("J3").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-13 19:19:59"

# Rebuild all four hyperlinks on this sheet (A2, I2, A3, I3) in one pass.
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/62288978-c74f-438f-83ff-b02031d3c663.md", "", "", "8cdfdab3-0d92-4118-9672-4add81af50ff.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/be4f88091d3164ed15249ce391f27d313d85cf8c/e2e/62288978-c74f-438f-83ff-b02031d3c663.md", "", "", "8cdfdab3-0d92-4118-9672-4add81af50ff.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/ad84fe27-7571-4181-b764-17442849d730.md", "", "", "ffffec629dde-426b-428a-a92c-5c9ae518e34f.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/be4f88091d3164ed15249ce391f27d313d85cf8c/e2e/ad84fe27-7571-4181-b764-17442849d730.md", "", "", "ffffec629dde-426b-428a-a92c-5c9ae518e34f.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2 (file #1)
$de.Range("G2").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.de-de.xlf"
$de.Range("H2").Value = "2016-08-13 19:19:38"
$de.Range("J2").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.de-de.xlf"
$de.Range("K2").Value = "2016-08-13 19:20:15"

# Row 3 (file #2)
$de.Range("G3").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.de-de.xlf"
$de.Range("H3").Value = "2016-08-13 19:19:38"
$de.Range("J3").Value = "8cdfdab3-0d92-4118-9672-4add81af50ff.2d09a8d6e74d9db0fb9e419fe3646ff7b9a913c8.de-de.xlf"
$de.Range("K3").Value = "2016-08-13 19:20:15"

# Rebuild all four hyperlinks on this sheet (A2, I2, A3, I3) in one pass.
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/62288978-c74f-438f-83ff-b02031d3c663.md", "", "", "8cdfdab3-0d92-4118-9672-4add81af50ff.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c537e82c908fcfedc93db84fa4fb6089e3f89077/e2e/62288978-c74f-438f-83ff-b02031d3c663.md", "", "", "8cdfdab3-0d92-4118-9672-4add81af50ff.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/5e48484403cb391a585572f481efeef969036719/e2e/ad84fe27-7571-4181-b764-17442849d730.md", "", "", "ffffec629dde-426b-428a-a92c-5c9ae518e34f.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/c537e82c908fcfedc93db84fa4fb6089e3f89077/e2e/ad84fe27-7571-4181-b764-17442849d730.md", "", "", "ffffec629dde-426b-428a-a92c-5c9ae518e34f.md") | Out-Null
